$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 8) to the sheet, mirroring the existing rows.
$ws.Range("A8").Value = 42612.890752314815
$ws.Range("B8").Value = -6
$ws.Range("C8").Value = 52
$ws.Range("D8").Value = 46
$ws.Range("E8").Value = 22
$ws.Range("F8").Value = 77
$ws.Range("G8").Value = 11623
$ws.Range("H8").Value = 9933
$ws.Range("I8").Value = 1566
$ws.Range("J8").Value = 151
$ws.Range("K8").Value = 134
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 7
$ws.Range("N8").Value = "Noun"
